$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 485.77777
$ws.Range("I80").Value = 377.1111
$ws.Range("K80").Value = 1131.3333
$ws.Range("M80").Value = -133.3333

$ws.Range("H83").Value = 485.77777
$ws.Range("I83").Value = 377.1111
$ws.Range("K83").Value = 3393.9999
$ws.Range("M83").Value = 1598.0001

$ws.Range("H92").Value = 13889508
$ws.Range("I92").Value = 2315328.2
$ws.Range("K92").Value = 2315328.2
$ws.Range("M92").Value = -2314080.2

$ws.Range("H116").Value = 9914.643
$ws.Range("I116").Value = 16000.714
$ws.Range("K116").Value = 16000.714
$ws.Range("M116").Value = -12558.714

$ws.Range("H118").Value = 823.7273
$ws.Range("I118").Value = 716.6667
$ws.Range("J118").Value = 1305.5
$ws.Range("K118").Value = 2150.0001
$ws.Range("L118").Value = 3916.5
$ws.Range("M118").Value = -493.0001000000002
$ws.Range("N118").Value = -7230.5

$ws.Range("H129").Value = 1098.5714
$ws.Range("I129").Value = 965.6667
$ws.Range("K129").Value = 2897.0001
$ws.Range("M129").Value = 2102.9999

$ws.Range("H133").Value = 43750
$ws.Range("J133").Value = 43750
$ws.Range("L133").Value = 43750
$ws.Range("N133").Value = -53870

$ws.Range("H135").Value = 2029.1333
$ws.Range("I135").Value = 1486.8529
$ws.Range("J135").Value = 3705.2727
$ws.Range("K135").Value = 13381.6761
$ws.Range("L135").Value = 33347.4543
$ws.Range("M135").Value = -10846.6761
$ws.Range("N135").Value = -38417.4543

$ws.Range("H137").Value = 1592.6052
$ws.Range("I137").Value = 1212.5385
$ws.Range("J137").Value = 2416.0833
$ws.Range("K137").Value = 3637.6155
$ws.Range("L137").Value = 7248.249899999999
$ws.Range("M137").Value = -1087.6155
$ws.Range("N137").Value = -12348.2499

$ws.Range("H138").Value = 2227.2935
$ws.Range("I138").Value = 950.6905
$ws.Range("J138").Value = 3299.64
$ws.Range("K138").Value = 2852.0715
$ws.Range("L138").Value = 9898.92
$ws.Range("M138").Value = 2287.9285
$ws.Range("N138").Value = -20178.92

$ws.Range("H141").Value = 1544.1136
$ws.Range("I141").Value = 1130.7188
$ws.Range("J141").Value = 2646.5
$ws.Range("K141").Value = 3392.1564
$ws.Range("L141").Value = 7939.5
$ws.Range("M141").Value = 1787.8436
$ws.Range("N141").Value = -18299.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2992.95
$ws.Range("I32").Value = 2424.1978
$ws.Range("J32").Value = 6486.7144
$ws.Range("K32").Value = 2424.1978
$ws.Range("L32").Value = 6486.7144
$ws.Range("M32").Value = -2137.1978
$ws.Range("N32").Value = -7060.7144

$ws.Range("H61").Value = 4984.516
$ws.Range("I61").Value = 6934.1055
$ws.Range("J61").Value = 1897.6666
$ws.Range("K61").Value = 6934.1055
$ws.Range("L61").Value = 1897.6666
$ws.Range("M61").Value = -6722.1055
$ws.Range("N61").Value = -2321.6666

$ws.Range("H74").Value = 998.25
$ws.Range("I74").Value = 844.65216
$ws.Range("J74").Value = 1390.7778
$ws.Range("K74").Value = 844.65216
$ws.Range("L74").Value = 1390.7778
$ws.Range("M74").Value = 29.34784000000002
$ws.Range("N74").Value = -3138.7778

$ws.Range("H77").Value = 998.25
$ws.Range("I77").Value = 844.65216
$ws.Range("J77").Value = 1390.7778
$ws.Range("K77").Value = 4223.2608
$ws.Range("L77").Value = 6953.889
$ws.Range("M77").Value = 144.7392
$ws.Range("N77").Value = -15689.889

$ws.Range("H136").Value = 4984.516
$ws.Range("I136").Value = 6934.1055
$ws.Range("J136").Value = 1897.6666
$ws.Range("K136").Value = 20802.3165
$ws.Range("L136").Value = 5692.9998
$ws.Range("M136").Value = -18252.3165
$ws.Range("N136").Value = -10792.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 42300
$ws.Range("J63").Value = 42300
$ws.Range("L63").Value = 42300
$ws.Range("N63").Value = -43672

$ws.Range("H66").Value = 42300
$ws.Range("J66").Value = 42300
$ws.Range("L66").Value = 126900
$ws.Range("N66").Value = -133764

$ws.Range("H94").Value = 1309.0476
$ws.Range("I94").Value = 459.16666
$ws.Range("K94").Value = 459.16666
$ws.Range("M94").Value = -8.166659999999979

$ws.Range("H126").Value = 30596.54
$ws.Range("J126").Value = 30596.54
$ws.Range("L126").Value = 30596.54
$ws.Range("N126").Value = -40476.54

$ws.Range("H134").Value = 4945.973
$ws.Range("I134").Value = 6828
$ws.Range("J134").Value = 2731.8235
$ws.Range("K134").Value = 20484
$ws.Range("L134").Value = 8195.470499999999
$ws.Range("M134").Value = -17949
$ws.Range("N134").Value = -13265.4705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 265140.97
$ws.Range("I31").Value = 1605.4762
$ws.Range("J31").Value = 916228.7
$ws.Range("K31").Value = 1605.4762
$ws.Range("L31").Value = 916228.7
$ws.Range("M31").Value = -1310.4762
$ws.Range("N31").Value = -916818.7

$ws.Range("H34").Value = 265140.97
$ws.Range("I34").Value = 1605.4762
$ws.Range("J34").Value = 916228.7
$ws.Range("K34").Value = 1605.4762
$ws.Range("L34").Value = 916228.7
$ws.Range("M34").Value = -1403.4762
$ws.Range("N34").Value = -916632.7

$ws.Range("H58").Value = 1760.921
$ws.Range("I58").Value = 1089.35
$ws.Range("J58").Value = 2507.111
$ws.Range("K58").Value = 1089.35
$ws.Range("L58").Value = 2507.111
$ws.Range("M58").Value = -886.3499999999999
$ws.Range("N58").Value = -2913.111

$ws.Range("H63").Value = 42000
$ws.Range("J63").Value = 42000
$ws.Range("L63").Value = 42000
$ws.Range("N63").Value = -43372

$ws.Range("H66").Value = 42000
$ws.Range("J66").Value = 42000
$ws.Range("L66").Value = 126000
$ws.Range("N66").Value = -132864

$ws.Range("H114").Value = 27166.666
$ws.Range("J114").Value = 27166.666
$ws.Range("L114").Value = 27166.666
$ws.Range("N114").Value = -35844.666

$ws.Range("H132").Value = 2539.9302
$ws.Range("I132").Value = 2160.3928
$ws.Range("J132").Value = 3248.4
$ws.Range("K132").Value = 6481.178400000001
$ws.Range("L132").Value = 9745.200000000001
$ws.Range("M132").Value = -3951.178400000001
$ws.Range("N132").Value = -14805.2

$ws.Range("H134").Value = 2225.7874
$ws.Range("I134").Value = 2659.8125
$ws.Range("J134").Value = 1299.8667
$ws.Range("K134").Value = 7979.4375
$ws.Range("L134").Value = 3899.6001
$ws.Range("M134").Value = -5444.4375
$ws.Range("N134").Value = -8969.6001

$ws.Range("H136").Value = 1760.921
$ws.Range("I136").Value = 1089.35
$ws.Range("J136").Value = 2507.111
$ws.Range("K136").Value = 3268.05
$ws.Range("L136").Value = 7521.333
$ws.Range("M136").Value = -718.0499999999997
$ws.Range("N136").Value = -12621.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 112.375
$ws.Range("I11").Value = 112.375
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 337.125
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -197.125

$ws.Range("H117").Value = 18527706
$ws.Range("I117").Value = 20325.8
$ws.Range("J117").Value = 25645930
$ws.Range("K117").Value = 60977.39999999999
$ws.Range("L117").Value = 76937790
$ws.Range("M117").Value = -57535.39999999999
$ws.Range("N117").Value = -76944674

$ws.Range("H118").Value = 872.5
$ws.Range("I118").Value = 930
$ws.Range("K118").Value = 2790
$ws.Range("M118").Value = -1547

$ws.Range("H121").Value = 943.4186
$ws.Range("I121").Value = 733.3333
$ws.Range("J121").Value = 959.175
$ws.Range("K121").Value = 2199.9999
$ws.Range("L121").Value = 2877.525
$ws.Range("M121").Value = -889.9998999999998
$ws.Range("N121").Value = -5497.525

$ws.Range("H129").Value = 30304572
$ws.Range("J129").Value = 2138.8333
$ws.Range("L129").Value = 6416.499899999999
$ws.Range("N129").Value = -16416.4999

$ws.Range("H131").Value = 3449197.2
$ws.Range("I131").Value = 16667013
$ws.Range("J131").Value = 1071.3043
$ws.Range("K131").Value = 50001039
$ws.Range("L131").Value = 3213.9129
$ws.Range("M131").Value = -49995999
$ws.Range("N131").Value = -13293.9129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 21800
$ws.Range("J63").Value = 21800
$ws.Range("L63").Value = 21800
$ws.Range("N63").Value = -23172

$ws.Range("H66").Value = 21800
$ws.Range("J66").Value = 21800
$ws.Range("L66").Value = 65400
$ws.Range("N66").Value = -72264

$ws.Range("H97").Value = 910
$ws.Range("I97").Value = 910
$ws.Range("K97").Value = 910
$ws.Range("M97").Value = -414

$ws.Range("H132").Value = 1887.2587
$ws.Range("I132").Value = 1463.1143
$ws.Range("J132").Value = 2532.6956
$ws.Range("K132").Value = 4389.3429
$ws.Range("L132").Value = 7598.0868
$ws.Range("M132").Value = -1859.3429
$ws.Range("N132").Value = -12658.0868

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 42200
$ws.Range("J64").Value = 42200
$ws.Range("L64").Value = 42200
$ws.Range("N64").Value = -42650

$ws.Range("H67").Value = 42200
$ws.Range("J67").Value = 42200
$ws.Range("L67").Value = 42200
$ws.Range("N67").Value = -43760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 37987.25
$ws.Range("J63").Value = 37987.25
$ws.Range("L63").Value = 37987.25
$ws.Range("N63").Value = -39235.25

$ws.Range("H66").Value = 37987.25
$ws.Range("J66").Value = 37987.25
$ws.Range("L66").Value = 113961.75
$ws.Range("N66").Value = -120201.75

$ws.Range("H126").Value = 953.7778
$ws.Range("I126").Value = 572.3333
$ws.Range("J126").Value = 1716.6666
$ws.Range("K126").Value = 1716.9999
$ws.Range("L126").Value = 5149.9998
$ws.Range("M126").Value = 753.0001
$ws.Range("N126").Value = -10089.9998

$ws.Range("H136").Value = 7465050
$ws.Range("I136").Value = 2383.311
$ws.Range("J136").Value = 22729596
$ws.Range("K136").Value = 7149.933000000001
$ws.Range("L136").Value = 68188788
$ws.Range("M136").Value = -4599.933000000001
$ws.Range("N136").Value = -68193888

$ws.Range("H138").Value = 39414.5
$ws.Range("J138").Value = 39414.5
$ws.Range("L138").Value = 39414.5
$ws.Range("N138").Value = -49694.5
